# Apply scheduling updates to the "Export.xlsx" schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value swaps (style changes accompany some of these) ---

# C8: "OFF" -> "9:45AM-4PM" ; style must change from the OFF fill (s=6) to the
# plain working-shift style (s=4), same as used by sibling cells e.g. B5.
$ws.Range("B5").Copy()
$ws.Range("C8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C8").Value = "9:45AM-4PM"

# C12: "9:45AM-4PM" -> "9:30AM-4PM" ; style (s=4) is unchanged.
$ws.Range("C12").Value = "9:30AM-4PM"

# C15: "10AM-5PM" -> "OFF" ; style must change to the OFF fill (s=6), same as
# used by sibling cells e.g. B16.
$ws.Range("B16").Copy()
$ws.Range("C15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C15").Value = "OFF"

# C16: "OFF" -> "10AM-5PM" ; style must change to the plain working-shift
# style (s=4), same as used by sibling cells e.g. B15.
$ws.Range("B15").Copy()
$ws.Range("C16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C16").Value = "10AM-5PM"

# --- "Unassigned Shifts" block: move the Bartender shift from C28 up into
#     C27 (the old C27 Lifeguard entry already lives on in D27), then remove
#     the now-empty C28 cell entirely. ---

$ws.Range("C27").Value = "Bartender,`n10AM-4PM"
$ws.Range("C28").Clear()
